$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 2.15
$ws.Range("H2").Value2 = 3.3
$ws.Range("I2").Value2 = 3.5
$ws.Range("N2").Value2 = 2.08
$ws.Range("O2").Value2 = 1.73
$ws.Range("P2").Value2 = 1.44
$ws.Range("Q2").Value2 = 2.63
$ws.Range("R2").Value2 = 1.83
$ws.Range("S2").Value2 = 1.83
$ws.Range("T2").Value2 = 7
$ws.Range("U2").Value2 = 10
$ws.Range("V2").Value2 = 9
$ws.Range("W2").Value2 = 19
$ws.Range("AD2").Value2 = 301
$ws.Range("AF2").Value2 = 17
$ws.Range("AH2").Value2 = 41
$ws.Range("AI2").Value2 = 29
$ws.Range("AJ2").Value2 = 41

# Row 3
$ws.Range("N3").Value2 = 1.73
$ws.Range("O3").Value2 = 2.08

# Row 14
$ws.Range("G14").Value2 = 3.65
$ws.Range("I14").Value2 = 1.88
$ws.Range("N14").Value2 = 2.12
$ws.Range("Q14").Value2 = 2.37
$ws.Range("T14").Value2 = 7.6
$ws.Range("V14").Value2 = 11
$ws.Range("X14").Value2 = 30
$ws.Range("AC14").Value2 = 75
$ws.Range("AE14").Value2 = 5.1
$ws.Range("AF14").Value2 = 6.8
$ws.Range("AH14").Value2 = 12.5
$ws.Range("AI14").Value2 = 13.5

# Row 15
$ws.Range("G15").Value2 = 3.3
$ws.Range("H15").Value2 = 2.9
$ws.Range("I15").Value2 = 2.4
$ws.Range("L15").Value2 = 1.44
$ws.Range("M15").Value2 = 2.63
$ws.Range("N15").Value2 = 2.4
$ws.Range("O15").Value2 = 1.53
$ws.Range("P15").Value2 = 1.53
$ws.Range("Q15").Value2 = 2.38
$ws.Range("R15").Value2 = 2.05
$ws.Range("S15").Value2 = 1.7
$ws.Range("T15").Value2 = 8
$ws.Range("Z15").Value2 = 6.5
$ws.Range("AA15").Value2 = 5.5
$ws.Range("AB15").Value2 = 17
$ws.Range("AC15").Value2 = 67
$ws.Range("AD15").Value2 = 451
$ws.Range("AG15").Value2 = 10
$ws.Range("AH15").Value2 = 23
$ws.Range("AI15").Value2 = 23
$ws.Range("AJ15").Value2 = 41

# Row 17
$ws.Range("G17").Value2 = 2.1
$ws.Range("H17").Value2 = 3.15
$ws.Range("I17").Value2 = 3.3
$ws.Range("J17").Value2 = 1.07
$ws.Range("K17").Value2 = 6.9
$ws.Range("L17").Value2 = 1.33
$ws.Range("M17").Value2 = 3.05
$ws.Range("N17").Value2 = 1.95
$ws.Range("O17").Value2 = 1.75
$ws.Range("P17").Value2 = 1.45
$ws.Range("Q17").Value2 = 2.57
$ws.Range("R17").Value2 = 1.75
$ws.Range("S17").Value2 = 1.95
$ws.Range("T17").Value2 = 7.1
$ws.Range("U17").Value2 = 10
$ws.Range("V17").Value2 = 8.75
$ws.Range("W17").Value2 = 20
$ws.Range("X17").Value2 = 18
$ws.Range("Y17").Value2 = 29
$ws.Range("Z17").Value2 = 6.9
$ws.Range("AA17").Value2 = 6.2
$ws.Range("AB17").Value2 = 13.5
$ws.Range("AC17").Value2 = 65
$ws.Range("AD17").Value2 = 500
$ws.Range("AE17").Value2 = 10
$ws.Range("AF17").Value2 = 18
$ws.Range("AG17").Value2 = 11.25
$ws.Range("AH17").Value2 = 45
$ws.Range("AI17").Value2 = 29
$ws.Range("AJ17").Value2 = 35

# Row 18
$ws.Range("G18").Value2 = 1.87
$ws.Range("H18").Value2 = 3.3
$ws.Range("I18").Value2 = 3.85
$ws.Range("J18").Value2 = 1.07
$ws.Range("K18").Value2 = 6.8
$ws.Range("L18").Value2 = 1.34
$ws.Range("M18").Value2 = 3
$ws.Range("N18").Value2 = 2
$ws.Range("O18").Value2 = 1.72
$ws.Range("P18").Value2 = 1.45
$ws.Range("Q18").Value2 = 2.55
$ws.Range("R18").Value2 = 1.85
$ws.Range("S18").Value2 = 1.85
$ws.Range("T18").Value2 = 6.5
$ws.Range("V18").Value2 = 8.5
$ws.Range("W18").Value2 = 16
$ws.Range("X18").Value2 = 15.5
$ws.Range("Y18").Value2 = 29
$ws.Range("Z18").Value2 = 6.8
$ws.Range("AA18").Value2 = 6.5
$ws.Range("AB18").Value2 = 15.5
$ws.Range("AC18").Value2 = 75
$ws.Range("AD18").Value2 = 700
$ws.Range("AE18").Value2 = 10.75
$ws.Range("AF18").Value2 = 21
$ws.Range("AH18").Value2 = 60
$ws.Range("AI18").Value2 = 37
$ws.Range("AJ18").Value2 = 45

# Row 19
$ws.Range("G19").Value2 = 1.7
$ws.Range("H19").Value2 = 3.5
$ws.Range("I19").Value2 = 4.4
$ws.Range("J19").Value2 = 1.06
$ws.Range("K19").Value2 = 7.2
$ws.Range("L19").Value2 = 1.31
$ws.Range("M19").Value2 = 3.15
$ws.Range("N19").Value2 = 1.93
$ws.Range("O19").Value2 = 1.78
$ws.Range("P19").Value2 = 1.42
$ws.Range("Q19").Value2 = 2.65
$ws.Range("R19").Value2 = 1.87
$ws.Range("S19").Value2 = 1.83
$ws.Range("T19").Value2 = 6.5
$ws.Range("U19").Value2 = 7.7
$ws.Range("V19").Value2 = 8.25
$ws.Range("W19").Value2 = 13
$ws.Range("X19").Value2 = 14
$ws.Range("Y19").Value2 = 28
$ws.Range("Z19").Value2 = 7.2
$ws.Range("AA19").Value2 = 6.9
$ws.Range("AB19").Value2 = 16
$ws.Range("AC19").Value2 = 80
$ws.Range("AD19").Value2 = 700
$ws.Range("AE19").Value2 = 12
$ws.Range("AF19").Value2 = 26
$ws.Range("AH19").Value2 = 75
$ws.Range("AJ19").Value2 = 50

# Row 20
$ws.Range("I20").Value2 = 5.4

# Row 24
$ws.Range("T24").Value2 = 7.2
$ws.Range("U24").Value2 = 9.5
$ws.Range("AA24").Value2 = 5.7
$ws.Range("AD24").Value2 = 250
$ws.Range("AF24").Value2 = 12.5

# Row 31
$ws.Range("G31").Value2 = 1.3
$ws.Range("H31").Value2 = 5.5
$ws.Range("I31").Value2 = 6
$ws.Range("K31").Value2 = 15
$ws.Range("L31").Value2 = 1.08
$ws.Range("M31").Value2 = 7
$ws.Range("N31").Value2 = 1.3
$ws.Range("O31").Value2 = 3.4
$ws.Range("R31").Value2 = 1.53
$ws.Range("S31").Value2 = 2.38
$ws.Range("U31").Value2 = 10
$ws.Range("V31").Value2 = 9.5
$ws.Range("W31").Value2 = 11
$ws.Range("Y31").Value2 = 19
$ws.Range("AA31").Value2 = 13
$ws.Range("AB31").Value2 = 17
$ws.Range("AC31").Value2 = 41
$ws.Range("AD31").Value2 = 101
$ws.Range("AF31").Value2 = 41
$ws.Range("AG31").Value2 = 21
$ws.Range("AH31").Value2 = 67
$ws.Range("AI31").Value2 = 41
$ws.Range("AJ31").Value2 = 34

# Row 32
$ws.Range("G32").Value2 = 2.1
$ws.Range("H32").Value2 = 3.3
$ws.Range("I32").Value2 = 3.3
$ws.Range("N32").Value2 = 1.9
$ws.Range("O32").Value2 = 1.9
$ws.Range("T32").Value2 = 8
$ws.Range("U32").Value2 = 11
$ws.Range("W32").Value2 = 19
$ws.Range("X32").Value2 = 17
$ws.Range("Y32").Value2 = 26
$ws.Range("AE32").Value2 = 11
$ws.Range("AF32").Value2 = 17
$ws.Range("AH32").Value2 = 34
$ws.Range("AI32").Value2 = 26

# Row 33
$ws.Range("G33").Value2 = 9
$ws.Range("H33").Value2 = 7
$ws.Range("I33").Value2 = 1.2
$ws.Range("K33").Value2 = 17
$ws.Range("L33").Value2 = 1.07
$ws.Range("M33").Value2 = 7.5
$ws.Range("N33").Value2 = 1.29
$ws.Range("O33").Value2 = 3.5
$ws.Range("R33").Value2 = 1.67
$ws.Range("S33").Value2 = 2.1
$ws.Range("T33").Value2 = 34
$ws.Range("U33").Value2 = 51
$ws.Range("V33").Value2 = 26
$ws.Range("W33").Value2 = 101
$ws.Range("X33").Value2 = 51
$ws.Range("AA33").Value2 = 15
$ws.Range("AB33").Value2 = 21
$ws.Range("AC33").Value2 = 51
$ws.Range("AD33").Value2 = 151
$ws.Range("AF33").Value2 = 9
$ws.Range("AG33").Value2 = 10
$ws.Range("AH33").Value2 = 9
$ws.Range("AJ33").Value2 = 21

# Row 34
$ws.Range("N34").Value2 = 2
$ws.Range("O34").Value2 = 1.8

# Row 35
$ws.Range("K35").Value2 = 15
$ws.Range("L35").Value2 = 1.18
$ws.Range("M35").Value2 = 4.5
$ws.Range("N35").Value2 = 1.62
$ws.Range("O35").Value2 = 2.25

# Row 36
$ws.Range("G36").Value2 = 2
$ws.Range("I36").Value2 = 3.5
$ws.Range("N36").Value2 = 1.73
$ws.Range("O36").Value2 = 2.08
$ws.Range("T36").Value2 = 8.5
$ws.Range("U36").Value2 = 10
$ws.Range("V36").Value2 = 8.5
$ws.Range("AF36").Value2 = 19
$ws.Range("AH36").Value2 = 41
$ws.Range("AI36").Value2 = 26

# Row 37
$ws.Range("G37").Value2 = 1.65
$ws.Range("H37").Value2 = 4.1
$ws.Range("N37").Value2 = 1.8
$ws.Range("O37").Value2 = 2
$ws.Range("P37").Value2 = 1.33
$ws.Range("Q37").Value2 = 3.25
$ws.Range("W37").Value2 = 12
$ws.Range("AF37").Value2 = 26
